$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.100696206092834
$ws.Range("B1").Value = 1.335281729698181
$ws.Range("C1").Value = 1.111676335334778
$ws.Range("D1").Value = 1.078254342079163
$ws.Range("E1").Value = 1.152105212211609
